$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 26.144619
$ws.Cells.Item(2, 8).Value = 78.433857
$ws.Cells.Item(2, 9).Value = 0.5211737020083955
$ws.Cells.Item(2, 10).Value = 0.5211737020083955
$ws.Cells.Item(2, 13).Value = 19.92674333333333
$ws.Cells.Item(2, 14).Value = 59.78023
$ws.Cells.Item(2, 15).Value = 0.3447897148135736
$ws.Cells.Item(2, 16).Value = 0.3447897148135735
$ws.Cells.Item(2, 17).Value = 520.9771123607901
$ws.Cells.Item(2, 18).Value = 4688.794011247111
$ws.Cells.Item(2, 19).Value = 0.1796953320838091
$ws.Cells.Item(2, 20).Value = 0.179695332083809

$ws.Cells.Item(3, 7).Value = 26.144619
$ws.Cells.Item(3, 8).Value = 78.433857
$ws.Cells.Item(3, 9).Value = 0.5211737020083955
$ws.Cells.Item(3, 10).Value = 0.5211737020083955
$ws.Cells.Item(3, 15).Value = 0.5793221821339875
$ws.Cells.Item(3, 16).Value = 0.5793221821339873
$ws.Cells.Item(3, 17).Value = 875.3555706785102
$ws.Cells.Item(3, 18).Value = 7878.20013610659
$ws.Cells.Item(3, 19).Value = 0.3019274863183522
$ws.Cells.Item(3, 20).Value = 0.3019274863183521

$ws.Cells.Item(4, 7).Value = 26.144619
$ws.Cells.Item(4, 8).Value = 78.433857
$ws.Cells.Item(4, 9).Value = 0.5211737020083955
$ws.Cells.Item(4, 10).Value = 0.5211737020083955
$ws.Cells.Item(4, 13).Value = 4.385869666666667
$ws.Cells.Item(4, 14).Value = 13.157609
$ws.Cells.Item(4, 15).Value = 0.07588810305243907
$ws.Cells.Item(4, 16).Value = 0.07588810305243905
$ws.Cells.Item(4, 17).Value = 114.666891418657
$ws.Cells.Item(4, 18).Value = 1032.002022767913
$ws.Cells.Item(4, 19).Value = 0.03955088360623429
$ws.Cells.Item(4, 20).Value = 0.03955088360623428

$ws.Cells.Item(5, 9).Value = 0.3571392594830743
$ws.Cells.Item(5, 10).Value = 0.3571392594830742
$ws.Cells.Item(5, 13).Value = 19.92674333333333
$ws.Cells.Item(5, 14).Value = 59.78023
$ws.Cells.Item(5, 15).Value = 0.3447897148135736
$ws.Cells.Item(5, 16).Value = 0.3447897148135735
$ws.Cells.Item(5, 17).Value = 357.0045445485001
$ws.Cells.Item(5, 18).Value = 3213.0409009365
$ws.Cells.Item(5, 19).Value = 0.1231379434259
$ws.Cells.Item(5, 20).Value = 0.1231379434259

$ws.Cells.Item(6, 9).Value = 0.3571392594830743
$ws.Cells.Item(6, 10).Value = 0.3571392594830742
$ws.Cells.Item(6, 15).Value = 0.5793221821339875
$ws.Cells.Item(6, 16).Value = 0.5793221821339873
$ws.Cells.Item(6, 17).Value = 599.8457694465001
$ws.Cells.Item(6, 19).Value = 0.206898695129451
$ws.Cells.Item(6, 20).Value = 0.2068986951294509

$ws.Cells.Item(7, 9).Value = 0.3571392594830743
$ws.Cells.Item(7, 10).Value = 0.3571392594830742
$ws.Cells.Item(7, 13).Value = 4.385869666666667
$ws.Cells.Item(7, 14).Value = 13.157609
$ws.Cells.Item(7, 15).Value = 0.07588810305243907
$ws.Cells.Item(7, 16).Value = 0.07588810305243905
$ws.Cells.Item(7, 17).Value = 78.57658306755002
$ws.Cells.Item(7, 18).Value = 707.18924760795
$ws.Cells.Item(7, 19).Value = 0.02710262092772332
$ws.Cells.Item(7, 20).Value = 0.02710262092772331

$ws.Cells.Item(8, 7).Value = 6.104416333333333
$ws.Cells.Item(8, 8).Value = 18.313249
$ws.Cells.Item(8, 9).Value = 0.1216870385085301
$ws.Cells.Item(8, 10).Value = 0.1216870385085301
$ws.Cells.Item(8, 13).Value = 19.92674333333333
$ws.Cells.Item(8, 14).Value = 59.78023
$ws.Cells.Item(8, 15).Value = 0.3447897148135736
$ws.Cells.Item(8, 16).Value = 0.3447897148135735
$ws.Cells.Item(8, 17).Value = 121.6411374741411
$ws.Cells.Item(8, 18).Value = 1094.77023726727
$ws.Cells.Item(8, 19).Value = 0.04195643930386445
$ws.Cells.Item(8, 20).Value = 0.04195643930386445

$ws.Cells.Item(9, 7).Value = 6.104416333333333
$ws.Cells.Item(9, 8).Value = 18.313249
$ws.Cells.Item(9, 9).Value = 0.1216870385085301
$ws.Cells.Item(9, 10).Value = 0.1216870385085301
$ws.Cells.Item(9, 15).Value = 0.5793221821339875
$ws.Cells.Item(9, 16).Value = 0.5793221821339873
$ws.Cells.Item(9, 17).Value = 204.38373353707
$ws.Cells.Item(9, 18).Value = 1839.45360183363
$ws.Cells.Item(9, 19).Value = 0.07049600068618425
$ws.Cells.Item(9, 20).Value = 0.07049600068618424

$ws.Cells.Item(10, 7).Value = 6.104416333333333
$ws.Cells.Item(10, 8).Value = 18.313249
$ws.Cells.Item(10, 9).Value = 0.1216870385085301
$ws.Cells.Item(10, 10).Value = 0.1216870385085301
$ws.Cells.Item(10, 13).Value = 4.385869666666667
$ws.Cells.Item(10, 14).Value = 13.157609
$ws.Cells.Item(10, 15).Value = 0.07588810305243907
$ws.Cells.Item(10, 16).Value = 0.07588810305243905
$ws.Cells.Item(10, 17).Value = 26.77317442907123
$ws.Cells.Item(10, 18).Value = 240.958569861641
$ws.Cells.Item(10, 19).Value = 0.009234598518481456
$ws.Cells.Item(10, 20).Value = 0.009234598518481456

